$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "38.189.96"
$ws.Range("E2").Value = "  +2.02%  "

$ws.Range("D3").Value = "2.056.68"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.58%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0823"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.73%  "

$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.00%  "

$ws.Range("D13").Value = "2.360.37"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.73%  "

$ws.Range("E15").Value = "  +2.53%  "

$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "2.055.54"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "38.135.51"
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").Value = "1.539.34"
$ws.Range("E40").Value = "  +4.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "

$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0931"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.77%  "

$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").Value = "2.248.83"
$ws.Range("E51").Value = "  +1.16%  "
